$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")
$wsStock = $wb.Worksheets.Item("股票")

# --- Row 1 ---------------------------------------------------------------
# The old row 1 simply duplicated row 2's data. It now becomes a proper
# header row (like every other sheet) and gains the same trailing
# metadata columns (G:M) the rest of the "normal" sheets carry.
$ws.Range("B1").Value = "bank"
$ws.Range("C1").Value = "deposit_type"
$ws.Range("D1").Value = "currency"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"
$ws.Range("G1").Value = "property_category"
$ws.Range("H1").Value = "category"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"

# "date" (I1) needs to land as a shared string rather than be
# auto-recognised as a real date, so pull the *value* only from a cell
# that already holds that exact literal (sheet "股票" I1 = "date").
$wsStock.Range("I1").Copy()
$ws.Range("I1").PasteSpecial(-4163) # xlPasteValues

# give the new header cells (G1:M1) the same bold/bordered look as the
# rest of row 1 (B1:F1 already carry style index 1)
$ws.Range("B1").Copy()
$ws.Range("G1:M1").PasteSpecial(-4122) # xlPasteFormats

# --- Rows 2-15 -------------------------------------------------------------
# Columns B:F already hold the right data (bank, deposit type, currency,
# owner, total). Only the new trailing columns G:M are missing and need
# to be filled in, mirroring every other "normal" sheet (property_category,
# category, date, legislator_name, legislator_id, source_file, index).
$firstRow = 2
$lastRow = 15
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Range("G$r").Value = "deposit"
    $ws.Range("H$r").Value = "normal"
    $ws.Range("J$r").Value = "賴士葆"
    $ws.Range("K$r").Value = 866
    $ws.Range("L$r").Value = "tmp9edb1"
    $ws.Range("M$r").Value = $r + 45   # same running index already used in column A

    # "2012-04-19" needs to stay a literal shared string (it already is one
    # used throughout the workbook) instead of becoming a date serial, so
    # copy the *value* from an existing cell holding that exact string.
    $wsStock.Range("J2").Copy()
    $ws.Range("I$r").PasteSpecial(-4163) # xlPasteValues
}

# Copy the data-row style (B2:F2 already uses style index 2) onto the new
# G:M cells for every data row.
$ws.Range("B2").Copy()
$ws.Range("G$firstRow`:M$lastRow").PasteSpecial(-4122) # xlPasteFormats
